# Insert a new weekly data row above the current row 66, shifting all
# subsequent rows down by one (old row 66 -> 67, ..., old row 216 -> 217),
# then populate the newly-inserted row 66 with the new price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new entire row before row 66 (pushes rows 66..216 down to 67..217)
$ws.Rows.Item(66).Insert()

# Fill in the new row 66 with the new record's data
$ws.Range("A66").Value = 8
$ws.Range("B66").Value = "Terminal La Palmera de La Serena"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44498
$ws.Range("D66").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E66").Value = 4
$ws.Range("F66").Value = 100114013
$ws.Range("G66").Value = "Zanahoria"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 700
$ws.Range("K66").Value = 6800
$ws.Range("L66").Value = 7000
$ws.Range("M66").Value = 6900
$ws.Range("N66").Value = "`$/saco 20 kilos"
$ws.Range("O66").Value = "Provincia del Elquí"
$ws.Range("P66").Value = 345
$ws.Range("Q66").Value = 20
$ws.Range("R66").Value = "Hortaliza"
